$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append case records (rows 703-729) that were previously missing because a bug
# held their conditions back from being written out.

# Row 703
$ws.Cells.Item(703, "A").Value = "21CRB01268"
$ws.Cells.Item(703, "B").Value = "Hemmeter"
$ws.Cells.Item(703, "C").Value = "POSSESSION DRUG PARAPHERNALIA"
$ws.Cells.Item(703, "D").Value = "2925.14(C)"
$ws.Cells.Item(703, "E").Value = "M4"
$ws.Cells.Item(703, "F").Value = "No Contest"
$ws.Cells.Item(703, "G").Value = "Guilty"
$cell = $ws.Cells.Item(703, "H")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"
$cell = $ws.Cells.Item(703, "I")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"

# Row 704
$ws.Cells.Item(704, "A").Value = "21CRB01268"
$ws.Cells.Item(704, "B").Value = "Hemmeter"
$ws.Cells.Item(704, "C").Value = "POSSESSION DRUG PARAPHERNALIA"
$ws.Cells.Item(704, "D").Value = "2925.14(C)"
$ws.Cells.Item(704, "E").Value = "M4"
$ws.Cells.Item(704, "F").Value = "No Contest"
$ws.Cells.Item(704, "G").Value = "Guilty"
$cell = $ws.Cells.Item(704, "H")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"
$cell = $ws.Cells.Item(704, "I")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"

# Row 705
$ws.Cells.Item(705, "A").Value = "03TRD13368"
$ws.Cells.Item(705, "B").Value = "Hemmeter"
$ws.Cells.Item(705, "C").Value = "SPEED REDUCED ZONE 3RD OR MORE"
$ws.Cells.Item(705, "D").Value = "4511.21C***"
$ws.Cells.Item(705, "E").Value = "M3"
$ws.Cells.Item(705, "F").Value = "No Contest"
$ws.Cells.Item(705, "G").Value = "Guilty"
$cell = $ws.Cells.Item(705, "H")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"
$cell = $ws.Cells.Item(705, "I")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"

# Row 706
$ws.Cells.Item(706, "A").Value = "03TRD13368"
$ws.Cells.Item(706, "B").Value = "Hemmeter"
$ws.Cells.Item(706, "C").Value = "SPEED REDUCED ZONE 3RD OR MORE"
$ws.Cells.Item(706, "D").Value = "4511.21C***"
$ws.Cells.Item(706, "E").Value = "M3"
$ws.Cells.Item(706, "F").Value = "No Contest"
$ws.Cells.Item(706, "G").Value = "Guilty"
$cell = $ws.Cells.Item(706, "H")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"
$cell = $ws.Cells.Item(706, "I")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"

# Row 707
$ws.Cells.Item(707, "A").Value = "21CRB01387"
$ws.Cells.Item(707, "B").Value = "Hemmeter"
$ws.Cells.Item(707, "C").Value = "SEXUAL IMPOSITION M1"
$ws.Cells.Item(707, "D").Value = "2907.06(A)(1)"
$ws.Cells.Item(707, "E").Value = "M1"
$ws.Cells.Item(707, "F").Value = "Guilty"
$ws.Cells.Item(707, "G").Value = "Guilty"
$cell = $ws.Cells.Item(707, "H")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"
$cell = $ws.Cells.Item(707, "I")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"

# Row 708
$ws.Cells.Item(708, "A").Value = "21CRB01268"
$ws.Cells.Item(708, "B").Value = "Bunner"
$ws.Cells.Item(708, "C").Value = "POSSESSION DRUG PARAPHERNALIA"
$ws.Cells.Item(708, "D").Value = "2925.14(C)"
$ws.Cells.Item(708, "E").Value = "M4"
$ws.Cells.Item(708, "F").Value = "No Contest"
$ws.Cells.Item(708, "G").Value = "Guilty"
$cell = $ws.Cells.Item(708, "H")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"
$cell = $ws.Cells.Item(708, "I")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"

# Row 709
$ws.Cells.Item(709, "A").Value = "21CRB01268"
$ws.Cells.Item(709, "B").Value = "Bunner"
$ws.Cells.Item(709, "C").Value = "POSSESSION DRUG PARAPHERNALIA"
$ws.Cells.Item(709, "D").Value = "2925.14(C)"
$ws.Cells.Item(709, "E").Value = "M4"
$ws.Cells.Item(709, "F").Value = "No Contest"
$ws.Cells.Item(709, "G").Value = "Guilty"
$cell = $ws.Cells.Item(709, "H")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"
$cell = $ws.Cells.Item(709, "I")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"

# Row 710
$ws.Cells.Item(710, "A").Value = "21CRB01268"
$ws.Cells.Item(710, "B").Value = "Bunner"
$ws.Cells.Item(710, "C").Value = "POSSESSION DRUG PARAPHERNALIA"
$ws.Cells.Item(710, "D").Value = "2925.14(C)"
$ws.Cells.Item(710, "E").Value = "M4"
$ws.Cells.Item(710, "F").Value = "No Contest"
$ws.Cells.Item(710, "G").Value = "Guilty"
$cell = $ws.Cells.Item(710, "H")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"
$cell = $ws.Cells.Item(710, "I")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"

# Row 711
$ws.Cells.Item(711, "A").Value = "21CRB01437"
$ws.Cells.Item(711, "B").Value = "Bunner"
$ws.Cells.Item(711, "C").Value = "POSSESSION OF MARIHUANA"
$ws.Cells.Item(711, "D").Value = "2925.11C3"
$ws.Cells.Item(711, "E").Value = "MM"
$ws.Cells.Item(711, "F").Value = "Guilty"
$ws.Cells.Item(711, "G").Value = "Guilty"
$cell = $ws.Cells.Item(711, "H")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"
$cell = $ws.Cells.Item(711, "I")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"

# Row 712
$ws.Cells.Item(712, "A").Value = "21CRB01437"
$ws.Cells.Item(712, "B").Value = "Bunner"
$ws.Cells.Item(712, "C").Value = "POSSESSION OF MARIHUANA"
$ws.Cells.Item(712, "D").Value = "2925.11C3"
$ws.Cells.Item(712, "E").Value = "MM"
$ws.Cells.Item(712, "F").Value = "Guilty"
$ws.Cells.Item(712, "G").Value = "Guilty"
$cell = $ws.Cells.Item(712, "H")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"
$cell = $ws.Cells.Item(712, "I")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"

# Row 713
$ws.Cells.Item(713, "A").Value = "21CRB01437"
$ws.Cells.Item(713, "B").Value = "Bunner"
$ws.Cells.Item(713, "C").Value = "POSSESSION OF MARIHUANA"
$ws.Cells.Item(713, "D").Value = "2925.11C3"
$ws.Cells.Item(713, "E").Value = "MM"
$ws.Cells.Item(713, "F").Value = "Guilty"
$ws.Cells.Item(713, "G").Value = "Guilty"
$cell = $ws.Cells.Item(713, "H")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"
$cell = $ws.Cells.Item(713, "I")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"

# Row 714
$ws.Cells.Item(714, "A").Value = "21CRB01437"
$ws.Cells.Item(714, "B").Value = "Bunner"
$ws.Cells.Item(714, "C").Value = "POSSESSION OF MARIHUANA"
$ws.Cells.Item(714, "D").Value = "2925.11C3"
$ws.Cells.Item(714, "E").Value = "MM"
$ws.Cells.Item(714, "F").Value = "Guilty"
$ws.Cells.Item(714, "G").Value = "Guilty"
$cell = $ws.Cells.Item(714, "H")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"
$cell = $ws.Cells.Item(714, "I")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"

# Row 715
$ws.Cells.Item(715, "A").Value = "21CRB01268"
$ws.Cells.Item(715, "B").Value = "Bunner"
$ws.Cells.Item(715, "C").Value = "POSSESSION DRUG PARAPHERNALIA"
$ws.Cells.Item(715, "D").Value = "2925.14(C)"
$ws.Cells.Item(715, "E").Value = "M4"
$ws.Cells.Item(715, "F").Value = "Guilty"
$ws.Cells.Item(715, "G").Value = "Guilty"
$cell = $ws.Cells.Item(715, "H")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"
$cell = $ws.Cells.Item(715, "I")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"

# Row 716
$ws.Cells.Item(716, "A").Value = "21CRB01268"
$ws.Cells.Item(716, "B").Value = "Hemmeter"
$ws.Cells.Item(716, "C").Value = "POSSESSION DRUG PARAPHERNALIA"
$ws.Cells.Item(716, "D").Value = "2925.14(C)"
$ws.Cells.Item(716, "E").Value = "M4"
$ws.Cells.Item(716, "F").Value = "No Contest"
$ws.Cells.Item(716, "G").Value = "Guilty"
$cell = $ws.Cells.Item(716, "H")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"
$cell = $ws.Cells.Item(716, "I")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"

# Row 717
$ws.Cells.Item(717, "A").Value = "21CRB01268"
$ws.Cells.Item(717, "B").Value = "Hemmeter"
$ws.Cells.Item(717, "C").Value = "POSSESSION DRUG PARAPHERNALIA"
$ws.Cells.Item(717, "D").Value = "2925.14(C)"
$ws.Cells.Item(717, "E").Value = "M4"
$ws.Cells.Item(717, "F").Value = "No Contest"
$ws.Cells.Item(717, "G").Value = "Guilty"
$cell = $ws.Cells.Item(717, "H")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"
$cell = $ws.Cells.Item(717, "I")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"

# Row 718
$ws.Cells.Item(718, "A").Value = "21CRB01268"
$ws.Cells.Item(718, "B").Value = "Hemmeter"
$ws.Cells.Item(718, "C").Value = "POSSESSION DRUG PARAPHERNALIA"
$ws.Cells.Item(718, "D").Value = "2925.14(C)"
$ws.Cells.Item(718, "E").Value = "M4"
$ws.Cells.Item(718, "F").Value = "No Contest"
$ws.Cells.Item(718, "G").Value = "Guilty"
$cell = $ws.Cells.Item(718, "H")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"
$cell = $ws.Cells.Item(718, "I")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"

# Row 719
$ws.Cells.Item(719, "A").Value = "21TRD09386"
$ws.Cells.Item(719, "B").Value = "Hemmeter"
$ws.Cells.Item(719, "C").Value = "DUS UCM"
$cell = $ws.Cells.Item(719, "D")
$cell.NumberFormat = "@"
$cell.Value = "4510.111"
$ws.Cells.Item(719, "E").Value = "UCM"
$ws.Cells.Item(719, "F").Value = "No Contest"
$ws.Cells.Item(719, "G").Value = "Guilty"
$cell = $ws.Cells.Item(719, "H")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"
$cell = $ws.Cells.Item(719, "I")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"

# Row 720
$ws.Cells.Item(720, "A").Value = "21TRD09386"
$ws.Cells.Item(720, "B").Value = "Hemmeter"
$ws.Cells.Item(720, "C").Value = "TAIL LIGHTS-REAR LICENSE PLATE"
$cell = $ws.Cells.Item(720, "D")
$cell.NumberFormat = "@"
$cell.Value = "4513.05"
$ws.Cells.Item(720, "E").Value = "MM"
$ws.Cells.Item(720, "F").Value = "No Contest"
$ws.Cells.Item(720, "G").Value = "Guilty"
$cell = $ws.Cells.Item(720, "H")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"
$cell = $ws.Cells.Item(720, "I")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"

# Row 721
$ws.Cells.Item(721, "A").Value = "21TRD09386"
$ws.Cells.Item(721, "B").Value = "Hemmeter"
$ws.Cells.Item(721, "C").Value = "DUS UCM"
$cell = $ws.Cells.Item(721, "D")
$cell.NumberFormat = "@"
$cell.Value = "4510.111"
$ws.Cells.Item(721, "E").Value = "UCM"
$ws.Cells.Item(721, "F").Value = "No Contest"
$ws.Cells.Item(721, "G").Value = "Guilty"
$cell = $ws.Cells.Item(721, "H")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"
$cell = $ws.Cells.Item(721, "I")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"

# Row 722
$ws.Cells.Item(722, "A").Value = "21TRD09386"
$ws.Cells.Item(722, "B").Value = "Hemmeter"
$ws.Cells.Item(722, "C").Value = "TAIL LIGHTS-REAR LICENSE PLATE"
$cell = $ws.Cells.Item(722, "D")
$cell.NumberFormat = "@"
$cell.Value = "4513.05"
$ws.Cells.Item(722, "E").Value = "MM"
$ws.Cells.Item(722, "F").Value = "No Contest"
$ws.Cells.Item(722, "G").Value = "Guilty"
$cell = $ws.Cells.Item(722, "H")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"
$cell = $ws.Cells.Item(722, "I")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"

# Row 723
$ws.Cells.Item(723, "A").Value = "21CRB01268"
$ws.Cells.Item(723, "B").Value = "Bunner"
$ws.Cells.Item(723, "C").Value = "POSSESSION DRUG PARAPHERNALIA"
$ws.Cells.Item(723, "D").Value = "2925.14(C)"
$ws.Cells.Item(723, "E").Value = "M4"
$ws.Cells.Item(723, "F").Value = "No Contest"
$ws.Cells.Item(723, "G").Value = "Guilty"
$cell = $ws.Cells.Item(723, "H")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"
$cell = $ws.Cells.Item(723, "I")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"

# Row 724
$ws.Cells.Item(724, "A").Value = "21TRD09246"
$ws.Cells.Item(724, "B").Value = "Bunner"
$ws.Cells.Item(724, "C").Value = "1ST SPEED IN 1 YR >70MPH"
$ws.Cells.Item(724, "D").Value = "4511.21D4"
$ws.Cells.Item(724, "E").Value = "No Data"
$ws.Cells.Item(724, "F").Value = "Guilty"
$ws.Cells.Item(724, "G").Value = "Guilty"
$cell = $ws.Cells.Item(724, "H")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"
$cell = $ws.Cells.Item(724, "I")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"

# Row 725
$ws.Cells.Item(725, "A").Value = "21TRD09246"
$ws.Cells.Item(725, "B").Value = "Bunner"
$ws.Cells.Item(725, "C").Value = "1ST SPEED IN 1 YR >70MPH"
$ws.Cells.Item(725, "D").Value = "4511.21D4"
$ws.Cells.Item(725, "E").Value = "No Data"
$ws.Cells.Item(725, "F").Value = "Guilty"
$ws.Cells.Item(725, "G").Value = "Guilty"
$cell = $ws.Cells.Item(725, "H")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"
$cell = $ws.Cells.Item(725, "I")
$cell.NumberFormat = "@"
$cell.Value = "`$ 0"

# Row 726
$ws.Cells.Item(726, "A").Value = "22CRB00136"
$ws.Cells.Item(726, "B").Value = "Bunner"
$ws.Cells.Item(726, "C").Value = "DOMESTIC VIOLENCE"
$ws.Cells.Item(726, "D").Value = "2919.25(A)"
$ws.Cells.Item(726, "E").Value = "No Data"
$ws.Cells.Item(726, "F").Value = "Not Guilty"

# Row 727
$ws.Cells.Item(727, "A").Value = "22CRB00136"
$ws.Cells.Item(727, "B").Value = "Bunner"
$ws.Cells.Item(727, "C").Value = "ASSAULT - M1"
$ws.Cells.Item(727, "D").Value = "2903.13(A)"
$ws.Cells.Item(727, "E").Value = "No Data"
$ws.Cells.Item(727, "F").Value = "Not Guilty"

# Row 728
$ws.Cells.Item(728, "A").Value = "22TRD00869"
$ws.Cells.Item(728, "B").Value = "Bunner"
$ws.Cells.Item(728, "C").Value = "DUS - DRIVING UNDER OVI SUSP"
$ws.Cells.Item(728, "D").Value = "4510.14A"
$ws.Cells.Item(728, "E").Value = "M1"
$ws.Cells.Item(728, "F").Value = "Not Guilty"

# Row 729
$ws.Cells.Item(729, "A").Value = "22TRD00869"
$ws.Cells.Item(729, "B").Value = "Bunner"
$ws.Cells.Item(729, "C").Value = "1ST SPEED 1 YR REDUCED ZONE"
$cell = $ws.Cells.Item(729, "D")
$cell.NumberFormat = "@"
$cell.Value = "4511.21"
$ws.Cells.Item(729, "E").Value = "MM"
$ws.Cells.Item(729, "F").Value = "Not Guilty"
